$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "'1.86%"
$ws.Range("D3").Value = "'27.37"
$ws.Range("E3").Value = "'1.64%"
$ws.Range("D4").Value = "'4.713"
$ws.Range("E4").Value = "'-1.40%"
$ws.Range("D5").Value = "'0.06092"
$ws.Range("E5").Value = "'3.02%"
$ws.Range("D6").Value = "'6.682"
$ws.Range("E6").Value = "'1.20%"
$ws.Range("D7").Value = "'0.8481"
$ws.Range("E7").Value = "'-0.35%"
$ws.Range("D8").Value = "'0.9263"
$ws.Range("E8").Value = "'0.38%"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'1.74%"
$ws.Range("D10").Value = "'0.04680"
$ws.Range("E10").Value = "'12.68%"
$ws.Range("D11").Value = "'0.07100"
$ws.Range("E11").Value = "'1.35%"
$ws.Range("D12").Value = "'0.03086"
$ws.Range("E12").Value = "'1.15%"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("D14").Value = "'0.001532"
$ws.Range("E14").Value = "'0.38%"
$ws.Range("D15").Value = "'0.0006084"
$ws.Range("E15").Value = "'0.28%"
$ws.Range("D16").Value = "'0.006093"
$ws.Range("E16").Value = "'1.24%"
$ws.Range("E17").Value = "'-0.61%"
$ws.Range("D18").Value = "'3.148"
$ws.Range("E18").Value = "'-0.53%"
$ws.Range("D20").Value = "'0.3110"
$ws.Range("E21").Value = "'0.20%"
$ws.Range("D22").Value = "'4.083"
$ws.Range("E22").Value = "'4.43%"
$ws.Range("D23").Value = "'0.04245"
$ws.Range("E23").Value = "'-0.32%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.22%"
$ws.Range("D25").Value = "'0.003789"
$ws.Range("E25").Value = "'-11.38%"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E27").Value = "'3.44%"
$ws.Range("D40").Value = "'0.03873"
$ws.Range("E40").Value = "'2.30%"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("E41").Value = "'1.43%"
$ws.Range("D42").Value = "'0.004082"
$ws.Range("E42").Value = "'-34.76%"
$ws.Range("D43").Value = "'0.01625"
$ws.Range("E43").Value = "'15.23%"
$ws.Range("E44").Value = "'0.80%"
$ws.Range("D45").Value = "'0.00005153"
$ws.Range("E45").Value = "'-3.57%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("D47").Value = "'0.1371"
$ws.Range("E47").Value = "'-38.67%"
$ws.Range("E48").Value = "'23.81%"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("E50").Value = "'0.07%"
